$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.758369565010071
$ws.Range("B1").Value = 2.138585329055786
$ws.Range("C1").Value = 2.257286787033081
$ws.Range("D1").Value = 2.66700005531311
$ws.Range("E1").Value = 2.923641681671143
